$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (shifts F:O -> G:P, updates formulas automatically)
$ws.Columns("F").Insert()

# Update header label text
$ws.Range("A39").Value = "Layer ID"

# Add new "Type 3" header in the freshly inserted column F (row 39 is the header row)
$ws.Range("F39").Value = "Type 3"

# Fix density unit label text (old L39 -> now M39 after the column insert)
$ws.Range("M39").Value = "Snow Density [g/cm3]"

# Restore the active cell selection
$ws.Range("M56").Select() | Out-Null
